$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of quarterly data: 01-07-2021 ("Q3 2021")
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "01-07-2021"

$values = @(43243, 9032, 34211, 10, 10, 0, 13268, 0, 13268, 697, 53, 644, 54856, 1, 54855, 204, 0, 204, 112279, 9097, 103181)

$col = 2
foreach ($v in $values) {
    $ws.Cells.Item(36, $col).Value = $v
    $col++
}
